# Add a new "year2" column (with placeholder "-" values) to the base-year
# Fill/Update-table worksheets, enabling the 2100-horizon runs (P16).
#
# Each of these worksheets holds a flat table with header row 1
# (scenario, attribute, process, commodity, lim_type, time_slice, year,
# commodity_group, currency, stage, sow, IE, National) in columns A:M.
# We append a new column N: header "year2" in row 1, and "-" (the
# workbook's standard "not applicable" placeholder) for every data row.

$wb = $excel.ActiveWorkbook

$sheetNames = @(
    "BY-RSD-WH_AF",
    "BY-RSD-SH_AF",
    "BY-RSD-EFF",
    "BY-RSD-DW",
    "BY-RSD-PF",
    "BY-RSD-CD",
    "BY-RSD-LT",
    "BY-RSD-CW",
    "BY-RSD-CK",
    "BY-RSD-RF",
    "BY-RSD-OE"
)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $lastRow = $ws.UsedRange.Rows.Count

    $ws.Cells.Item(1, 14).Value = "year2"

    for ($r = 2; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 14).Value = "-"
    }
}
